$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths (A, B, C)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 17.5
$ws.Columns.Item(2).ColumnWidth = 155.66666666666666
$ws.Columns.Item(3).ColumnWidth = 11.833333333333332

# ---------------------------------------------------------------------------
# Cell values - header / letterhead block (rows 1-4)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = 43269
$ws.Range("B1").Value = "Journal de bord"
$ws.Range("B2").Value = "Jeremy Comelli"

# ---------------------------------------------------------------------------
# Cell values - table header (row 5) + existing log entries (rows 6-8)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Date"
$ws.Range("B5").Value = "Entrée"

$ws.Range("A6").Value = 17.05
$ws.Range("B6").Value = "Le développement de l'interface graphique prend plus de temps que prévu, cependant la partie ""Chargeur de fichiers"" va finalement être abandonnée, car le processus est moins complexe que prévu. "

$ws.Range("A7").Value = 23.05
$ws.Range("B7").Value = "Après entretien avec M. Ithurbide, il a été décidé que la méthode de prendre un screenshot n'était effectivement pas optimisée, ou portable (celle-ci utilisait notamment un offset de coordonées hardodé, prévu pour ignorer spécifiquement les bordures de fenêtres windows 7). À la place, une conversion de l'image en matrice numpy sera utilisée"

$ws.Range("A8").Value = 24.05
$ws.Range("B8").Value = "J'ai passé la matinée sur un problème passablement frustrant. En essayant de modifier le chargeur d'image pour le faire créer des objets, j'ai oublié de convertir Image.fromarray(image) en ImageTk.PhotoImage. Problème résolu, mais du temps a été perdu inutilement sur un problème facile. Cependant, avec la nouvelle architecture orientée objet, il me sera plus facile de récupérer une partie de l'image."

# ---------------------------------------------------------------------------
# Cell values - new log entries (rows 9-10)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = 24.05
$ws.Range("B9").Value = "Note à moi-même: Il aurait été bien plus simple d'empêcher le retournement de la selection_box,  par exemple en swappant le côté sélectionné quand la souris passe au dessus, car il est complexe de transformer des coordonnées négatives, et une boite de sélection basée sur des valeurs width et height."

$ws.Range("A10").Value = 24.05
$ws.Range("B10").Value = "Après discussions avec M. Ithurbide, je me suis rendu compte que j'ai failli partir sur une feature qui ne se trouvait pas dans le cahier des charges (le collage d'images, qui se fera avec programme tiers, comme photoshop). Il a aussi été décidé de développer le recadrage, et d'implémenter le processing d'images à travers des matrices numpy."

# ---------------------------------------------------------------------------
# Blank styled rows below the table (rows 11-33, column B)
# ---------------------------------------------------------------------------
$ws.Range("B11:B12").WrapText = $true
$ws.Range("B11:B12").VerticalAlignment = -4108
$ws.Range("B13:B33").WrapText = $true

# ---------------------------------------------------------------------------
# Merge the letterhead date cell
# ---------------------------------------------------------------------------
$ws.Range("A1:A2").Merge()

# ---------------------------------------------------------------------------
# Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 29.25
$ws.Rows.Item(6).RowHeight = 31.5
$ws.Rows.Item(7).RowHeight = 47.25
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 31.5
$ws.Rows.Item(10).RowHeight = 47.25

# ---------------------------------------------------------------------------
# Date number format (A1) + header/date column look
# ---------------------------------------------------------------------------
$ws.Range("A1").NumberFormat = "mm-dd-yy"
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108

$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108

$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108

$ws.Range("B2").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Thin box border around the full table (header row 5 + data rows 6-10, col A+B)
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A5:B10")
$tableRange.Borders.Item(7).LineStyle = 1
$tableRange.Borders.Item(7).Weight = 2
$tableRange.Borders.Item(8).LineStyle = 1
$tableRange.Borders.Item(8).Weight = 2
$tableRange.Borders.Item(9).LineStyle = 1
$tableRange.Borders.Item(9).Weight = 2
$tableRange.Borders.Item(10).LineStyle = 1
$tableRange.Borders.Item(10).Weight = 2
$tableRange.Borders.Item(11).LineStyle = 1
$tableRange.Borders.Item(11).Weight = 2
$tableRange.Borders.Item(12).LineStyle = 1
$tableRange.Borders.Item(12).Weight = 2

# ---------------------------------------------------------------------------
# Medium outer box border around the letterhead (A1:B2 block)
# ---------------------------------------------------------------------------
$ws.Range("A1").Borders.Item(7).LineStyle = 1
$ws.Range("A1").Borders.Item(7).Weight = -4138
$ws.Range("A1").Borders.Item(8).LineStyle = 1
$ws.Range("A1").Borders.Item(8).Weight = -4138

$ws.Range("B1").Borders.Item(10).LineStyle = 1
$ws.Range("B1").Borders.Item(10).Weight = -4138
$ws.Range("B1").Borders.Item(8).LineStyle = 1
$ws.Range("B1").Borders.Item(8).Weight = -4138

$ws.Range("A2").Borders.Item(7).LineStyle = 1
$ws.Range("A2").Borders.Item(7).Weight = -4138

$ws.Range("B2").Borders.Item(10).LineStyle = 1
$ws.Range("B2").Borders.Item(10).Weight = -4138

# ---------------------------------------------------------------------------
# Thin inner box border around the little A3:B4 spacer block
# ---------------------------------------------------------------------------
$ws.Range("A3").Borders.Item(7).LineStyle = 1
$ws.Range("A3").Borders.Item(7).Weight = 2
$ws.Range("A3").Borders.Item(8).LineStyle = 1
$ws.Range("A3").Borders.Item(8).Weight = 2

$ws.Range("B3").Borders.Item(10).LineStyle = 1
$ws.Range("B3").Borders.Item(10).Weight = 2
$ws.Range("B3").Borders.Item(8).LineStyle = 1
$ws.Range("B3").Borders.Item(8).Weight = 2

$ws.Range("A4").Borders.Item(7).LineStyle = 1
$ws.Range("A4").Borders.Item(7).Weight = 2
$ws.Range("A4").Borders.Item(9).LineStyle = 1
$ws.Range("A4").Borders.Item(9).Weight = 2

$ws.Range("B4").Borders.Item(10).LineStyle = 1
$ws.Range("B4").Borders.Item(10).Weight = 2
$ws.Range("B4").Borders.Item(9).LineStyle = 1
$ws.Range("B4").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------------
# Bigger font + vertical centering for the entry texts
# ---------------------------------------------------------------------------
$ws.Range("B6:B10").Font.Size = 12
$ws.Range("B8:B10").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Print area / page setup
# ---------------------------------------------------------------------------
$ws.PageSetup.PrintArea = '$A$1:$C$28'
$ws.PageSetup.Orientation = 2
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.LeftMargin = 51.02362204724410
$ws.PageSetup.RightMargin = 150.23622047244098
$ws.PageSetup.TopMargin = 53.85826771653544
$ws.PageSetup.BottomMargin = 53.85826771653544
$ws.PageSetup.HeaderMargin = 22.677165354330707
$ws.PageSetup.FooterMargin = 22.677165354330707
$ws.PageSetup.Zoom = 77

# ---------------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------------
$ws.Range("B11").Select()
